# Update countries & provincias Spain
# - Reorder two pairs of countries whose updated "Casos totales" (col B)
#   ranking moved them past their former neighbours (Finlandia past Haiti,
#   Eslovenia past Siria/Gambia).
# - Refresh the numeric COVID columns (B..H) for the rows whose figures
#   changed in this data pull.
# - Bump the "Datos actualizados" timestamp string.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Country reordering -----------------------------------------------
# Gabon, Haiti, Finlandia, Zimbabue  ->  Gabon, Finlandia, Haiti, Zimbabue
$ws.Cells.Item(103, 1).Value = "Finlandia"
$ws.Cells.Item(104, 1).Value = "Haiti"

# Jamaica, Siria, Gambia, Eslovenia, Lituania -> Jamaica, Eslovenia, Siria, Gambia, Lituania
$ws.Cells.Item(129, 1).Value = "Eslovenia"
$ws.Cells.Item(130, 1).Value = "Siria"
$ws.Cells.Item(131, 1).Value = "Gambia"

# --- Numeric refresh (Casos totales, Nuevos casos, Casos activos,
#     Recuperados, Casos criticos, Muertes hoy, Muertes) --------------

# Row 18 - Banglades
$ws.Cells.Item(18, 2).Value = 331078
$ws.Cells.Item(18, 3).Value = 1827
$ws.Cells.Item(18, 4).Value = 230804
$ws.Cells.Item(18, 5).Value = 95681
$ws.Cells.Item(18, 7).Value = 41
$ws.Cells.Item(18, 8).Value = 4593

# Row 24 - Alemania
$ws.Cells.Item(24, 2).Value = 254957
$ws.Cells.Item(24, 3).Value = 1
$ws.Cells.Item(24, 5).Value = 14948

# Row 25 - Filipinas
$ws.Cells.Item(25, 2).Value = 245143
$ws.Cells.Item(25, 3).Value = 3176
$ws.Cells.Item(25, 4).Value = 185543
$ws.Cells.Item(25, 5).Value = 55614
$ws.Cells.Item(25, 7).Value = 70
$ws.Cells.Item(25, 8).Value = 3986

# Row 26 - Indonesia
$ws.Cells.Item(26, 2).Value = 203342
$ws.Cells.Item(26, 3).Value = 3307
$ws.Cells.Item(26, 4).Value = 145200
$ws.Cells.Item(26, 5).Value = 49806
$ws.Cells.Item(26, 7).Value = 106
$ws.Cells.Item(26, 8).Value = 8336

# Row 40 - Oman
$ws.Cells.Item(40, 2).Value = 87939
$ws.Cells.Item(40, 3).Value = 349
$ws.Cells.Item(40, 4).Value = 83115
$ws.Cells.Item(40, 5).Value = 4073
$ws.Cells.Item(40, 7).Value = 9
$ws.Cells.Item(40, 8).Value = 751

# Row 49 - Polonia
$ws.Cells.Item(49, 4).Value = 57135
$ws.Cells.Item(49, 5).Value = 12665

# Row 70 - Austria
$ws.Cells.Item(70, 2).Value = 30583
$ws.Cells.Item(70, 3).Value = 502
$ws.Cells.Item(70, 4).Value = 25764
$ws.Cells.Item(70, 5).Value = 4072

# Row 97 - Malasia
$ws.Cells.Item(97, 2).Value = 9583
$ws.Cells.Item(97, 3).Value = 24
$ws.Cells.Item(97, 4).Value = 9143
$ws.Cells.Item(97, 5).Value = 312

# Row 103 - Finlandia (new figures)
$ws.Cells.Item(103, 2).Value = 8430
$ws.Cells.Item(103, 3).Value = 93
$ws.Cells.Item(103, 4).Value = 7350
$ws.Cells.Item(103, 5).Value = 744
$ws.Cells.Item(103, 8).Value = 336

# Row 104 - Haiti (keeps its previous figures, now one row down)
$ws.Cells.Item(104, 2).Value = 8376
$ws.Cells.Item(104, 4).Value = 5991
$ws.Cells.Item(104, 5).Value = 2171
$ws.Cells.Item(104, 8).Value = 214

# Row 114 - Hong Kong
$ws.Cells.Item(114, 2).Value = 4902
$ws.Cells.Item(114, 3).Value = 6
$ws.Cells.Item(114, 4).Value = 4557
$ws.Cells.Item(114, 5).Value = 246

# Row 124 - Uganda
$ws.Cells.Item(124, 2).Value = 4101
$ws.Cells.Item(124, 3).Value = 201
$ws.Cells.Item(124, 4).Value = 1876
$ws.Cells.Item(124, 5).Value = 2179

# Row 129 - Eslovenia (new figures)
$ws.Cells.Item(129, 2).Value = 3312
$ws.Cells.Item(129, 3).Value = 79
$ws.Cells.Item(129, 4).Value = 2587
$ws.Cells.Item(129, 5).Value = 590
$ws.Cells.Item(129, 8).Value = 135

# Row 130 - Siria (keeps its previous figures, now one row down)
$ws.Cells.Item(130, 2).Value = 3289
$ws.Cells.Item(130, 4).Value = 760
$ws.Cells.Item(130, 5).Value = 2389
$ws.Cells.Item(130, 8).Value = 140

# Row 131 - Gambia (keeps its previous figures, now one row down)
$ws.Cells.Item(131, 2).Value = 3275
$ws.Cells.Item(131, 4).Value = 1424
$ws.Cells.Item(131, 5).Value = 1752
$ws.Cells.Item(131, 8).Value = 99

# Row 133 - Sri Lanka
$ws.Cells.Item(133, 4).Value = 2946
$ws.Cells.Item(133, 5).Value = 182

# Row 158 - Letonia
$ws.Cells.Item(158, 2).Value = 1443
$ws.Cells.Item(158, 3).Value = 11
$ws.Cells.Item(158, 4).Value = 1234
$ws.Cells.Item(158, 5).Value = 174

# --- Timestamp ----------------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 9 de Septiembre de 2020 a las 11:54"
